$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.324665069580078
$ws.Range("B1").Value = 1.420799851417542
$ws.Range("C1").Value = 1.634160041809082
$ws.Range("D1").Value = 2.687262296676636
$ws.Range("E1").Value = 4.584642887115479
